$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.318.88'
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").Value = '2.309.13'
$ws.Range("E3").Value = '  -2.13%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.81'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.53%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.610'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.58'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0916'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.31'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.89%  '
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.979'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.65'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.24%  '
$ws.Range("D16").Value = '2.658.80'
$ws.Range("E16").Value = '  -2.19%  '
$ws.Range("D17").Value = '2.311.03'
$ws.Range("E17").Value = '  -2.13%  '
$ws.Range("D18").Value = '42.202.39'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -5.45%  '
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '74.85'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.47'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -6.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '259.40'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.31'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.37'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -7.28%  '
$ws.Range("E26").Value = '  +0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.04'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.88%  '
$ws.Range("E28").Value = '  +3.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '22.89'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.80'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '165.40'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.87%  '
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.93'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.87'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.25%  '
$ws.Range("E35").Value = '  +11.77%  '
$ws.Range("E36").Value = '  -2.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.58'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0355'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("E39").Value = '  -5.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.63'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -4.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '72.73'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.84'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +7.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.47'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.18%  '
$ws.Range("E44").Value = '  -3.31%  '
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.34'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.05'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -6.32%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.05'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.35'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.47%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.03'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.66%  '
$ws.Range("D51").Value = '1.574.88'
$ws.Range("E51").Value = '  +0.73%  '
